$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.175.82"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "2.252.81"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'307.85"
$ws.Range("E5").Value = "  -4.33%  "
$ws.Range("D6").Value = "'98.69"
$ws.Range("E6").Value = "  -2.02%  "
$ws.Range("D7").Value = "'0.576"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("D10").Value = "'35.64"
$ws.Range("E10").Value = "  -3.25%  "
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "'7.32"
$ws.Range("E12").Value = "  -5.01%  "
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("D14").Value = "2.595.52"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "2.254.22"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "'0.839"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "'13.81"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "44.010.56"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").Value = "'12.99"
$ws.Range("E19").Value = "  -3.78%  "
$ws.Range("D20").Value = "0.0₃0976"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").Value = "'6.32"
$ws.Range("E21").Value = "  -3.42%  "
$ws.Range("D22").Value = "'65.34"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "'242.07"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("E24").Value = "  -6.90%  "
$ws.Range("E25").Value = "  -7.47%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'10.11"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").Value = "'36.70"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("D30").Value = "'6.22"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").Value = "'20.17"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").Value = "'156.73"
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("E33").Value = "  +14.30%  "
$ws.Range("D34").Value = "'0.0825"
$ws.Range("E34").Value = "  -2.93%  "
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("E37").Value = "  -4.37%  "
$ws.Range("D38").Value = "'1.85"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").Value = "'15.47"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("E40").Value = "  -8.27%  "
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("E42").Value = "  -9.95%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "1.767.11"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("D45").Value = "'87.16"
$ws.Range("E45").Value = "  +6.23%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'16.15"
$ws.Range("E46").Value = "  +13.32%  "
$ws.Range("D47").Value = "'5.16"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.193"
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'101.35"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'55.88"
$ws.Range("E51").Value = "  -4.24%  "
